$wb = $excel.ActiveWorkbook

# Rename the "Include from mCSD Endpoint Ty" sheet to "Include #0"
$wsInclude = $wb.Worksheets.Item("Include from mCSD Endpoint Ty")
$wsInclude.Name = "Include #0"

# Metadata sheet updates
$ws = $wb.Worksheets.Item("Metadata")

# Version: 3.8.0 -> 3.9.0
$ws.Range("B3").Value = "3.9.0"

# Experimental value (was blank) -> "false", written as literal text (not boolean).
# Build it as a formula result, then paste only the value so the text type/style
# of the destination cell is preserved instead of Excel auto-coercing "false" to a Boolean.
$helper = $ws.Range("D1")
$helper.Formula = "=""false"""
$helper.Copy()
$ws.Range("B7").PasteSpecial(-4163)
$helper.Clear()

# Date
$ws.Range("B8").Value = "2024-12-02T17:05:26-06:00"

# Contact rows (A10:A12 already say "Contact"); update the Value column for each
$ws.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$ws.Range("B11").Value = "null (iti@ihe.net)"
$ws.Range("B12").Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"

# Jurisdiction
$ws.Range("B13").Value = "Global (Whole world)"
